$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Extend formatting to a new row 16 by copying row 15's existing format
$ws.Range("A15:B15").Copy()
$ws.Range("A16:B16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Rewrite rows 12-16 with the shifted-down content (a new "Jurisdiction" row
# is inserted after "Contact", pushing Description/Purpose/Copyright/Immutable down)
$ws.Range("A12").Value = "Jurisdiction"
$ws.Range("B12").Value = ""
$ws.Range("A13").Value = "Description"
$ws.Range("B13").Value = "RxNorm codes for Alemtuzumab"
$ws.Range("A14").Value = "Purpose"
$ws.Range("B14").Value = ""
$ws.Range("A15").Value = "Copyright"
$ws.Range("B15").Value = ""
$ws.Range("A16").Value = "Immutable"
$ws.Range("B16").Value = "BooleanType[null]"

# Apply the other scalar value updates
$ws.Range("B3").Value = "0.1.7"
$ws.Range("B6").Value = "draft"
$ws.Range("B8").Value = "2024-08-27T12:23:18-05:00"
$ws.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"
$ws.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"
